$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "Estimates".
# This also re-points the _xlnm._FilterDatabase defined names that
# reference the sheet by name.
$ws.Name = "Estimates"

# Row 65's label was a duplicate of row 64's "Min (P=99%)" text but the
# formulas in C65/D65 actually compute the *upper* bound (I57+3*C62), so
# correct the label to "Max (P=99%)".
$ws.Range("A65").Value = "Max (P=99%)"
